$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "-Inserir touch no slider principal;" gets split into three
# runs ("-" / "Inserir" / " touch no slider principal;") with proofing
# marks around "Inserir" (wavy-underline spell-check artifacts are not
# reproducible through automation, so we focus on the run split, which
# is the structurally meaningful part of the change).
# ---------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("-Inserir touch no slider principal;", $false, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if ($found1) {
    $base1 = $rng1.Start

    $seg1a = $d.Range($base1, $base1 + 1)        # "-"
    $seg1b = $d.Range($base1 + 1, $base1 + 8)    # "Inserir"

    # Nudging bold on/off forces the engine to break these characters
    # into their own runs without leaving any real formatting behind.
    $seg1a.Font.Bold = $true
    $seg1a.Font.Bold = $false

    $seg1b.Font.Bold = $true
    $seg1b.Font.Bold = $false
}

# ---------------------------------------------------------------------
# Edit 2: "-Inserir o botão de favoritar em productPage do mobile (está
# display:none) – productPage linha 291 css;" gets colored 00B0F0 and
# the "9" in "291" is corrected to "2" (-> "221"), which naturally
# splits the middle run into three pieces around the edited character.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("productPage linha 291 css", $false, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if ($found2) {
    # Fix the line-number typo first, while the run still carries the
    # paragraph's original (uncoloured) formatting.
    $digitPos = $rng2.Start + 19
    $digitRange = $d.Range($digitPos, $digitPos + 1)
    $digitRange.Text = "2"
}

$rng3 = $d.Content
$found3 = $rng3.Find.Execute("-Inserir o bot", $false, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if ($found3) {
    $pstart = $rng3.Start

    $targetColor = 15773696   # RGB(0, 176, 240) == 00B0F0

    $bounds = @(0, 54, 96, 97, 102, 103)
    for ($i = 0; $i -lt ($bounds.Length - 1); $i++) {
        $piece = $d.Range($pstart + $bounds[$i], $pstart + $bounds[$i + 1])
        $piece.Font.Color = $targetColor
    }

    # Also color the paragraph mark itself (pPr/rPr) by colouring the
    # whole paragraph range, which includes the end-of-paragraph mark.
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Start -eq $pstart) {
            $p.Range.Font.Color = $targetColor
        }
    }
}
